$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Copy the row-31 formatting down into row 32 first, so the new row
# inherits the same cell styles (date format on column A, etc.)
$ws.Range("A31:K31").Copy()
$ws.Range("A32:K32").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# Add new attendance row for 06-09-2023 (row 32)
$ws.Range("A32").Value = 45175
$ws.Range("B32").Value = "PRESENT"
$ws.Range("C32").Value = "PRESENT"
$ws.Range("D32").Value = "PRESENT"
$ws.Range("E32").Value = "PRESENT"
$ws.Range("F32").Value = "PRESENT"
$ws.Range("G32").Value = "PRESENT"
$ws.Range("H32").Value = "ABSENT"
$ws.Range("I32").Value = "ABSENT"
$ws.Range("J32").Value = "ABSENT"
$ws.Range("K32").Value = "ABSENT"

$ws.Range("K32").Select()
